$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.994
$ws.Range("E5").Value = 12.949
$ws.Range("E9").Value = 12.82
$ws.Range("E11").Value = 13.003
$ws.Range("C21").Value = -13.075
$ws.Range("E21").Value = 13.258
$ws.Range("C23").Value = -12.937
$ws.Range("C25").Value = -13.175
